# Reorders the data rows (2-16) of the active worksheet: each target row
# receives the full set of "moving" column values (D, K:T) that originally
# belonged to a different source row, per the mapping below. Columns
# A, B, C, E, F, G, H, I, J are identical across all rows and therefore
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> original (source) row number
$rowMap = @{
    2  = 16
    3  = 13
    4  = 9
    5  = 10
    6  = 14
    7  = 11
    8  = 5
    9  = 6
    10 = 4
    11 = 3
    12 = 7
    13 = 8
    14 = 2
    15 = 12
    16 = 15
}

# Columns that move together with each row (1-based column indices):
# D=4, K=11, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot all original values first, since rows are interdependent.
# Value2 is used (rather than Value) because it returns/accepts the raw
# underlying number/string without any COM date-wrapper quirks.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write permuted values back.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $srcVals[$c]
    }
}
